$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: paint Status (column L) formatting by copying from existing Passed/Failed exemplars ---
$ws.Range("L70").Copy()
$ws.Range("L72:L77").PasteSpecial(-4122)
$ws.Range("L71").Copy()
$ws.Range("L78:L79").PasteSpecial(-4122)
$ws.Range("L70").Copy()
$ws.Range("L80").PasteSpecial(-4122)
$ws.Range("L71").Copy()
$ws.Range("L81:L82").PasteSpecial(-4122)
$ws.Range("L70").Copy()
$ws.Range("L83:L94").PasteSpecial(-4122)

# --- Step 2: set cell values (AR / ER / Status text) ---
$ws.Range("H72").Value = 'When clicking on submit button, the data entered by user should be sent to the database.
And the user will be redirected to the home page'
$ws.Range("I72").Value = 'When clicking on submit button, the data entered by user should be sent to the database.
And the user will be redirected to the home page'
$ws.Range("L72").Value = 'Passed'
$ws.Range("I73").Value = 'The user has been added succesfully to the database'
$ws.Range("L73").Value = 'Passed'
$ws.Range("I74").Value = 'The error message displayed successfully'
$ws.Range("L74").Value = 'Passed'
$ws.Range("I75").Value = 'The error message displayed successfully'
$ws.Range("L75").Value = 'Passed'
$ws.Range("I76").Value = 'The error message displayed successfully'
$ws.Range("L76").Value = 'Passed'
$ws.Range("I77").Value = 'The error message displayed successfully'
$ws.Range("L77").Value = 'Passed'
$ws.Range("I78").Value = 'The system accepted the username and added it to the database'
$ws.Range("L78").Value = 'Failed'
$ws.Range("I79").Value = 'The system accepted the username and added it to the database'
$ws.Range("L79").Value = 'Failed'
$ws.Range("I80").Value = 'The username has been accepted successfully by the system and added to the database'
$ws.Range("L80").Value = 'Passed'
$ws.Range("I81").Value = 'The system refused to proceed while leaving the optional fields empty'
$ws.Range("L81").Value = 'Failed'
$ws.Range("I82").Value = 'The system accepted the username although it is already existing'
$ws.Range("L82").Value = 'Failed'
$ws.Range("I83").Value = 'The system generated an error message to tell the user that the passwords entered don''t match'
$ws.Range("L83").Value = 'Passed'
$ws.Range("I84").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L84").Value = 'Passed'
$ws.Range("I85").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L85").Value = 'Passed'
$ws.Range("I86").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L86").Value = 'Passed'
$ws.Range("I87").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L87").Value = 'Passed'
$ws.Range("I88").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L88").Value = 'Passed'
$ws.Range("I89").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L89").Value = 'Passed'
$ws.Range("I90").Value = 'The system  proceeded successfully with submitting the form'
$ws.Range("L90").Value = 'Passed'
$ws.Range("I91").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L91").Value = 'Passed'
$ws.Range("I92").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L92").Value = 'Passed'
$ws.Range("I93").Value = 'The system refused to proceed and displayed the error message'
$ws.Range("L93").Value = 'Passed'
$ws.Range("I94").Value = 'The system  proceeded successfully with submitting the form'
$ws.Range("L94").Value = 'Passed'

# --- Step 3: update view state (scroll position + active selection) to match final state ---
$ws.Application.Goto($ws.Range("D92"), $true)
$ws.Range("I94").Select()
